# Generate Report for Handback
# Adds a new handback row (file f7f39c49-5320-4293-8d9a-99ea6278c1ef.md) to
# the "Overview", "zh-cn" and "de-de" sheets, mirroring the existing
# "04645944-e760-4fe7-bb22-f3e54f62209d.md" ("Handed back: in sync with
# en-US") row that is already present as row 3 of each table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New file identity / timestamps used throughout the new row
# ---------------------------------------------------------------------
$fileId          = "f7f39c49-5320-4293-8d9a-99ea6278c1ef"
$fileName        = "$fileId.md"
$pathAndName     = "e2e\$fileId.md"
$status          = "Handed back: in sync with en-US"
$genDate         = "2016-08-28 18:43:19"

$zhXliff         = "$fileId.973232f1552fcf6ef0c11db2edc2b678a63ff6c4.zh-cn.xlf"
$zhHandoffDate   = "2016-08-28 18:43:15"
$zhHandbackDate  = "2016-08-28 18:43:31"

$deXliff         = "$fileId.973232f1552fcf6ef0c11db2edc2b678a63ff6c4.de-de.xlf"
$deHandoffDate   = "2016-08-28 18:43:19"
$deHandbackDate  = "2016-08-28 18:43:37"

$srcUrl          = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1e6c4d2aed7b226992f0fbb213dd057d5b4cfa1/e2e/$fileId.md"
$zhUrl           = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6c01f8fbe881cc9de98c11093a1a66f392410a3f/e2e/$fileId.md"
$deUrl           = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f69998a412979f525f813f96aa892ee57a1b0d5d/e2e/$fileId.md"

# ---------------------------------------------------------------------
# Sheet "Overview" -> new row 4 (mirrors row 3)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $fileName
$wsOverview.Range("B4").Value = $pathAndName
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = $status
$wsOverview.Range("F4").Value = $status
$wsOverview.Range("G4").Value = $genDate
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Range("B4").Font.Underline = $true
$wsOverview.Range("B4").Font.Color = 6749952

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $srcUrl, "", "", $pathAndName) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" -> new row 4 (mirrors row 3)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A4").Value = $fileName
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = $status
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = $zhXliff
$wsZh.Range("H4").Value = $zhHandoffDate
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = $fileName
$wsZh.Range("J4").Value = $zhXliff
$wsZh.Range("K4").Value = $zhHandbackDate
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Range("A4").Font.Underline = $true
$wsZh.Range("A4").Font.Color = 6749952
$wsZh.Range("I4").Font.Underline = $true
$wsZh.Range("I4").Font.Color = 6749952

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $srcUrl, "", "", $fileName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), $zhUrl, "", "", $fileName) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de" -> new row 4 (mirrors row 3)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A4").Value = $fileName
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = $status
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = $deXliff
$wsDe.Range("H4").Value = $genDate
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = $fileName
$wsDe.Range("J4").Value = $deXliff
$wsDe.Range("K4").Value = $deHandbackDate
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Range("A4").Font.Underline = $true
$wsDe.Range("A4").Font.Color = 6749952
$wsDe.Range("I4").Font.Underline = $true
$wsDe.Range("I4").Font.Color = 6749952

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $srcUrl, "", "", $fileName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), $deUrl, "", "", $fileName) | Out-Null

Write-Host "Handback row added for $fileName"
